$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.598.56'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +0.61%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.922.99'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  -0.11%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("D4").ClearFormats()
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '247.35'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +2.78%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.000'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -0.01%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4743'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -0.12%  '
$ws.Range("E8").Value = '  +1.69%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06830'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +4.08%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '105.57'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -1.02%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '18.41'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -3.03%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.933.69'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +0.74%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.07701'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +1.17%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.358'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +4.58%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6704'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +2.52%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '291.01'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -3.68%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '30.622.93'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +0.68%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000007626'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +1.95%  '
$ws.Range("E19").Value = '  -0.05%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.96'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +0.24%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.564'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +5.14%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.179.30'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +0.41%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.0000'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +0.01%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.467'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +3.52%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.515'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +3.47%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '167.57'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +0.29%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.97'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +4.39%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.125'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +5.49%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.1072'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -3.16%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.403'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +3.57%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.178'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +2.49%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.053'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +3.69%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05015'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +0.77%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7337'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -0.76%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.147'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +0.26%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02072'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +7.14%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.9999'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +0.04%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.723'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -0.95%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.683'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -0.56%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '111.74'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +4.59%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.042'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -0.60%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.8732'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -0.31%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.4411'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +6.92%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.904'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +2.11%  '
$ws.Range("E45").Value = '  +0.00%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '67.83'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -2.83%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '7.295'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +0.78%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.349'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +0.72%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '48.33'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +15.34%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.1242'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +3.74%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '35.05'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +0.71%  '
